$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 for columns D, L, M, N, O, P, R, S
# (all other columns remain identical between the two rows)

$row2Vals = @{
    "D" = $ws.Range("D2").Value2
    "L" = $ws.Range("L2").Value2
    "M" = $ws.Range("M2").Value2
    "N" = $ws.Range("N2").Value2
    "O" = $ws.Range("O2").Value2
    "P" = $ws.Range("P2").Value2
    "R" = $ws.Range("R2").Value2
    "S" = $ws.Range("S2").Value2
}

$row3Vals = @{
    "D" = $ws.Range("D3").Value2
    "L" = $ws.Range("L3").Value2
    "M" = $ws.Range("M3").Value2
    "N" = $ws.Range("N3").Value2
    "O" = $ws.Range("O3").Value2
    "P" = $ws.Range("P3").Value2
    "R" = $ws.Range("R3").Value2
    "S" = $ws.Range("S3").Value2
}

foreach ($col in $row2Vals.Keys) {
    $ws.Range("$col`2").Value2 = $row3Vals[$col]
    $ws.Range("$col`3").Value2 = $row2Vals[$col]
}
